# Auto-generated Excel COM-interop script to apply profit recalculation updates
$wb = $excel.ActiveWorkbook

# ALC row 18
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 825.5
$ws.Range("I18").Value = 825.5
$ws.Range("K18").Value = 825.5
$ws.Range("M18").Value = -541.5

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3324.3333
$ws.Range("I74").Value = 3292.0833
$ws.Range("J74").Value = 3453.3333
$ws.Range("K74").Value = 3292.0833
$ws.Range("L74").Value = 3453.3333
$ws.Range("M74").Value = -2356.0833
$ws.Range("N74").Value = -5325.3333

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3324.3333
$ws.Range("I77").Value = 3292.0833
$ws.Range("J77").Value = 3453.3333
$ws.Range("K77").Value = 16460.4165
$ws.Range("L77").Value = 17266.6665
$ws.Range("M77").Value = -11780.4165
$ws.Range("N77").Value = -26626.6665

# ALC row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 309.03705
$ws.Range("I92").Value = 297
$ws.Range("K92").Value = 297
$ws.Range("M92").Value = 951

# ALC row 101
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 7575968.5
$ws.Range("I101").Value = 253.2
$ws.Range("J101").Value = 45454544
$ws.Range("K101").Value = 759.5999999999999
$ws.Range("L101").Value = 136363632
$ws.Range("M101").Value = 862.4000000000001
$ws.Range("N101").Value = -136366876

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 6414324
$ws.Range("I132").Value = 7249714
$ws.Range("J132").Value = 9667
$ws.Range("K132").Value = 21749142
$ws.Range("L132").Value = 29001
$ws.Range("M132").Value = -21746612
$ws.Range("N132").Value = -34061

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1184.8983
$ws.Range("I137").Value = 964.12
$ws.Range("J137").Value = 2411.4443
$ws.Range("K137").Value = 2892.36
$ws.Range("L137").Value = 7234.3329
$ws.Range("M137").Value = -342.3600000000001
$ws.Range("N137").Value = -12334.3329

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1459.9333
$ws.Range("I138").Value = 623.9423
$ws.Range("J138").Value = 3350
$ws.Range("K138").Value = 1871.8269
$ws.Range("L138").Value = 10050
$ws.Range("M138").Value = 3268.1731
$ws.Range("N138").Value = -20330

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2372.9
$ws.Range("I32").Value = 2009.3684
$ws.Range("J32").Value = 9280
$ws.Range("K32").Value = 2009.3684
$ws.Range("L32").Value = 9280
$ws.Range("M32").Value = -1722.3684
$ws.Range("N32").Value = -9854

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 766.44446
$ws.Range("I61").Value = 753.61536
$ws.Range("J61").Value = 1100
$ws.Range("K61").Value = 753.61536
$ws.Range("L61").Value = 1100
$ws.Range("M61").Value = -541.61536
$ws.Range("N61").Value = -1524

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1240.1471
$ws.Range("I74").Value = 1232.2727
$ws.Range("J74").Value = 1500
$ws.Range("K74").Value = 1232.2727
$ws.Range("L74").Value = 1500
$ws.Range("M74").Value = -358.2727
$ws.Range("N74").Value = -3248

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1240.1471
$ws.Range("I77").Value = 1232.2727
$ws.Range("J77").Value = 1500
$ws.Range("K77").Value = 6161.363499999999
$ws.Range("L77").Value = 7500
$ws.Range("M77").Value = -1793.363499999999
$ws.Range("N77").Value = -16236

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1308546.5
$ws.Range("I132").Value = 1056.6666
$ws.Range("K132").Value = 3169.9998
$ws.Range("M132").Value = -639.9998

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 766.44446
$ws.Range("I136").Value = 753.61536
$ws.Range("J136").Value = 1100
$ws.Range("K136").Value = 2260.84608
$ws.Range("L136").Value = 3300
$ws.Range("M136").Value = 289.1539199999997
$ws.Range("N136").Value = -8400

# BSM row 80
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 4195.8125
$ws.Range("I80").Value = 812.5
$ws.Range("J80").Value = 5323.5835
$ws.Range("K80").Value = 812.5
$ws.Range("L80").Value = 5323.5835
$ws.Range("M80").Value = 185.5
$ws.Range("N80").Value = -7319.5835

# BSM row 83
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 4195.8125
$ws.Range("I83").Value = 812.5
$ws.Range("J83").Value = 5323.5835
$ws.Range("K83").Value = 4062.5
$ws.Range("L83").Value = 26617.9175
$ws.Range("M83").Value = 929.5
$ws.Range("N83").Value = -36601.9175

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2713694
$ws.Range("I134").Value = 1001.3143
$ws.Range("J134").Value = 18537736
$ws.Range("K134").Value = 3003.9429
$ws.Range("L134").Value = 55613208
$ws.Range("M134").Value = -468.9429
$ws.Range("N134").Value = -55618278

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 9524876
$ws.Range("I132").Value = 886.0357
$ws.Range("J132").Value = 47620836
$ws.Range("K132").Value = 2658.1071
$ws.Range("L132").Value = 142862508
$ws.Range("M132").Value = -128.1071000000002
$ws.Range("N132").Value = -142867568

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 6329848
$ws.Range("I134").Value = 833.5208
$ws.Range("J134").Value = 16129613
$ws.Range("K134").Value = 2500.5624
$ws.Range("L134").Value = 48388839
$ws.Range("M134").Value = 34.4376000000002
$ws.Range("N134").Value = -48393909

# CUL row 32
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 61728896
$ws.Range("I32").Value = 1000
$ws.Range("J32").Value = 123456790
$ws.Range("K32").Value = 3000
$ws.Range("L32").Value = 370370370
$ws.Range("N32").Value = -370370936
$ws.Range("M32").Value = -2717

# CUL row 126
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 166668500
$ws.Range("I126").Value = 166668500
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 500005500
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -500000560
$ws.Range("N126").ClearContents()

# LTW row 32
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

# LTW row 35
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 7369.9
$ws.Range("I35").Value = 3083
$ws.Range("J35").Value = 24517.5
$ws.Range("K35").Value = 3083
$ws.Range("L35").Value = 24517.5
$ws.Range("M35").Value = -2747
$ws.Range("N35").Value = -25189.5

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2452734.8
$ws.Range("I46").Value = 5209201.5
$ws.Range("J46").Value = 2542.2222
$ws.Range("K46").Value = 5209201.5
$ws.Range("L46").Value = 2542.2222
$ws.Range("M46").Value = -5209013.5
$ws.Range("N46").Value = -2918.2222

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1032.6857
$ws.Range("I61").Value = 950.4762
$ws.Range("K61").Value = 950.4762
$ws.Range("M61").Value = -748.4762

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1386
$ws.Range("I82").Value = 1362
$ws.Range("J82").Value = 1410
$ws.Range("K82").Value = 1362
$ws.Range("L82").Value = 1410
$ws.Range("M82").Value = -1001
$ws.Range("N82").Value = -2132

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1386
$ws.Range("I85").Value = 1362
$ws.Range("J85").Value = 1410
$ws.Range("K85").Value = 1362
$ws.Range("L85").Value = 1410
$ws.Range("M85").Value = -114
$ws.Range("N85").Value = -3906

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1032.6857
$ws.Range("I113").Value = 950.4762
$ws.Range("K113").Value = 950.4762
$ws.Range("M113").Value = 1219.5238

# LTW row 116
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5518.423
$ws.Range("I122").Value = 5973
$ws.Range("J122").Value = 2033.3334
$ws.Range("K122").Value = 17919
$ws.Range("L122").Value = 6100.0002
$ws.Range("M122").Value = -15469
$ws.Range("N122").Value = -11000.0002

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6498.7446
$ws.Range("I132").Value = 2059.2917
$ws.Range("J132").Value = 11131.218
$ws.Range("K132").Value = 6177.875100000001
$ws.Range("L132").Value = 33393.654
$ws.Range("M132").Value = -3647.875100000001
$ws.Range("N132").Value = -38453.654

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 8131762.5
$ws.Range("I136").Value = 9525611
$ws.Range("J136").Value = 980
$ws.Range("K136").Value = 28576833
$ws.Range("L136").Value = 2940
$ws.Range("M136").Value = -28574283
$ws.Range("N136").Value = -8040

# WVR row 40
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 4350
$ws.Range("J40").Value = 5000
$ws.Range("L40").Value = 5000
$ws.Range("N40").Value = -5298

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1401.8889
$ws.Range("I126").Value = 945.2857
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 2835.8571
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -365.8571000000002
$ws.Range("N126").Value = -13940

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 13159498
$ws.Range("I136").Value = 17857964
$ws.Range("K136").Value = 53573892
$ws.Range("M136").Value = -53571342
